$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 233
$ws.Range("B233").Value = 6078269
$ws.Range("F233").Value = "Universidad de Chile"
$ws.Range("G233").Value = "Nublense"
$ws.Range("H233").Value = 3
$ws.Range("J233").Value = "H"
$ws.Range("K233").Value = 1.85
$ws.Range("L233").Value = 3.4
$ws.Range("M233").Value = 4.333
$ws.Range("N233").Value = 1.8
$ws.Range("O233").Value = 3.6
$ws.Range("P233").Value = 4.5
$ws.Range("Q233").Value = -0.75
$ws.Range("R233").Value = 1.925
$ws.Range("S233").Value = 1.925
$ws.Range("T233").Value = 2.5
$ws.Range("U233").Value = 2.025
$ws.Range("V233").Value = 1.825
$ws.Range("W233").Value = 0.8
$ws.Range("Y233").Value = -1
$ws.Range("Z233").Value = 0.925
$ws.Range("AA233").Value = -1
$ws.Range("AB233").Value = 1.025
$ws.Range("AC233").Value = -1

# Row 234
$ws.Range("B234").Value = 6077768
$ws.Range("F234").Value = "Union La Calera"
$ws.Range("G234").Value = "Universidad Catolica"
$ws.Range("H234").Value = 0
$ws.Range("I234").Value = 3
$ws.Range("J234").Value = "A"
$ws.Range("K234").Value = 2.05
$ws.Range("L234").Value = 3.5
$ws.Range("N234").Value = 2.05
$ws.Range("O234").Value = 3.6
$ws.Range("P234").Value = 3.4
$ws.Range("Q234").Value = -0.25
$ws.Range("R234").Value = 1.8
$ws.Range("S234").Value = 2
$ws.Range("U234").Value = 1.975
$ws.Range("V234").Value = 1.825
$ws.Range("W234").Value = -1
$ws.Range("Y234").Value = 2.4
$ws.Range("Z234").Value = -1
$ws.Range("AA234").Value = 1
$ws.Range("AB234").Value = 0.4875
$ws.Range("AC234").Value = -0.5

# Row 236
$ws.Range("B236").Value = 6077499
$ws.Range("F236").Value = "Deportes Copiapo"
$ws.Range("G236").Value = "Everton de Vina"
$ws.Range("H236").Value = 2
$ws.Range("I236").Value = 0
$ws.Range("J236").Value = "H"
$ws.Range("K236").Value = 2.1
$ws.Range("L236").Value = 3.4
$ws.Range("N236").Value = 2.9
$ws.Range("O236").Value = 3.4
$ws.Range("P236").Value = 2.4
$ws.Range("Q236").Value = 0.25
$ws.Range("R236").Value = 1.775
$ws.Range("S236").Value = 2.1
$ws.Range("U236").Value = 1.85
$ws.Range("V236").Value = 2
$ws.Range("W236").Value = 1.9
$ws.Range("Y236").Value = -1
$ws.Range("Z236").Value = 0.7749999999999999
$ws.Range("AA236").Value = -1
$ws.Range("AB236").Value = -1
$ws.Range("AC236").Value = 1

# Row 237
$ws.Range("B237").Value = 6078268
$ws.Range("F237").Value = "OHiggins"
$ws.Range("G237").Value = "Palestino"
$ws.Range("H237").Value = 0
$ws.Range("J237").Value = "A"
$ws.Range("K237").Value = 3.1
$ws.Range("L237").Value = 3.3
$ws.Range("M237").Value = 2.3
$ws.Range("N237").Value = 2.9
$ws.Range("O237").Value = 3.4
$ws.Range("P237").Value = 2.375
$ws.Range("Q237").Value = 0.25
$ws.Range("R237").Value = 1.8
$ws.Range("S237").Value = 2
$ws.Range("T237").Value = 2.75
$ws.Range("U237").Value = 2
$ws.Range("V237").Value = 1.8
$ws.Range("W237").Value = -1
$ws.Range("Y237").Value = 1.375
$ws.Range("Z237").Value = -1
$ws.Range("AA237").Value = 1
$ws.Range("AB237").Value = -1
$ws.Range("AC237").Value = 0.8

# Row 247
$ws.Range("B247").Value = 7723523
$ws.Range("F247").Value = "Deportes Iquique"
$ws.Range("G247").Value = "Everton de Vina"
$ws.Range("K247").Value = 3.3
$ws.Range("L247").Value = 3.3
$ws.Range("M247").Value = 2.2
$ws.Range("N247").Value = 2.9
$ws.Range("O247").Value = 3.3
$ws.Range("P247").Value = 2.45
$ws.Range("Q247").Value = 0
$ws.Range("R247").Value = 2.1
$ws.Range("S247").Value = 1.775
$ws.Range("U247").Value = 2
$ws.Range("V247").Value = 1.85

# Row 248
$ws.Range("B248").Value = 7723522
$ws.Range("F248").Value = "Huachipato"
$ws.Range("G248").Value = "Union La Calera"
$ws.Range("K248").Value = 2.1
$ws.Range("L248").Value = 3.4
$ws.Range("M248").Value = 3.5
$ws.Range("N248").Value = 1.85
$ws.Range("O248").Value = 3.6
$ws.Range("P248").Value = 4.2
$ws.Range("Q248").Value = -0.5
$ws.Range("R248").Value = 1.925
$ws.Range("S248").Value = 1.925
$ws.Range("U248").Value = 1.9
$ws.Range("V248").Value = 1.95

# Row 249
$ws.Range("R249").Value = 2
$ws.Range("S249").Value = 1.85
